$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2 = "H") with Week 17 data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 291
$wsOff.Range("C2").Value = 202
$wsOff.Range("D2").Value = 57
$wsOff.Range("F2").Value = 5

# Update DEF sheet (row 2 = "H") with Week 17 data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 198
$wsDef.Range("C2").Value = 123
$wsDef.Range("D2").Value = 50
$wsDef.Range("E2").Value = 19
$wsDef.Range("F2").Value = 5
$wsDef.Range("G2").Value = 1
